$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 2-5: update the genetic-algorithm sample table (Numero/Binario/X/f(x)).
# Binario (C) and f(x) (E) columns hold numeric-looking *text* in the sheet
# (e.g. "001110", "12.9600"). Assigning the literal string directly would
# let Excel auto-coerce it to a number (dropping leading/trailing zeros), so
# instead a quoted-text formula is entered and immediately flattened to a
# value in place via Copy + PasteSpecial(values). That keeps the result a
# plain shared-string text cell without leaving behind any extra
# NumberFormat-driven cell style.
# ---------------------------------------------------------------------------
function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $escaped = $text.Replace('"', '""')
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}

$ws.Range("B2").Value = 14
Set-TextValue "C2" "001110"
$ws.Range("D2").Value = 3.84
Set-TextValue "E2" "14.7456"

$ws.Range("B3").Value = 10
Set-TextValue "C3" "001010"
$ws.Range("D3").Value = 3.6
Set-TextValue "E3" "12.9600"

$ws.Range("B4").Value = 32
Set-TextValue "C4" "100000"
$ws.Range("D4").Value = 4.92
Set-TextValue "E4" "24.2064"

$ws.Range("B5").Value = 6
Set-TextValue "C5" "000110"
$ws.Range("D5").Value = 3.36
Set-TextValue "E5" "11.2896"

# ---------------------------------------------------------------------------
# Row 7: replace the small ColA/ColB header with the full GA parameters
# header row (A7:O7). Copy the existing header style (from B7, which already
# carries the bold/border/centered style) onto the newly used cells first.
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("D7:O7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A7").Value = "initial_people"
$ws.Range("B7").Value = "max_people"
$ws.Range("C7").Value = "variable_a"
$ws.Range("D7").Value = "variable_b"
$ws.Range("E7").Value = "prob_of_crossing"
$ws.Range("F7").Value = "prob_mut_ind"
$ws.Range("G7").Value = "prob_mut_gen"
$ws.Range("H7").Value = "function_entry"
$ws.Range("I7").Value = "delta_x"
$ws.Range("J7").Value = "find_x_by"
$ws.Range("K7").Value = "iterator_entry"
$ws.Range("L7").Value = "range"
$ws.Range("M7").Value = "jump_numbers"
$ws.Range("N7").Value = "points_numbers"
$ws.Range("O7").Value = "bits_required"

# ---------------------------------------------------------------------------
# Row 8: new data row under the GA parameters header. Most cells are text
# (including numeric-looking ones like "4", "0.75", "2"); I8/L8/M8/N8/O8
# are real numbers.
# ---------------------------------------------------------------------------
Set-TextValue "A8" "4"
Set-TextValue "B8" "8"
Set-TextValue "C8" "3"
Set-TextValue "D8" "5"
Set-TextValue "E8" "0.75"
Set-TextValue "F8" "0.25"
Set-TextValue "G8" "0.35"
$ws.Range("H8").Value = "x**2"
$ws.Range("I8").Value = 0.06
$ws.Range("J8").Value = "Minimización"
Set-TextValue "K8" "2"
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 33.3333
$ws.Range("N8").Value = 34.3333
$ws.Range("O8").Value = 6

# ---------------------------------------------------------------------------
# Rows 9 and 10 previously held a small lookup table (B/C) that is gone now.
# ---------------------------------------------------------------------------
$ws.Range("A9:O10").ClearContents()
